# Quarterly symbol-list refresh (GitHub Actions bot) — updates Price (D) and
# Volume(1h) (E) text columns for the changed ticker rows. Values are kept as
# literal text (NumberFormat "@") so they round-trip exactly like the scraped
# strings ("289.96", "-4.27%", ...) instead of being reinterpreted as numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    'D2' = '289.96'
    'E2' = '-4.27%'
    'D3' = '30.87'
    'E3' = '-5.82%'
    'D4' = '4.939'
    'D5' = '0.07196'
    'E5' = '-8.09%'
    'D6' = '1.788'
    'E6' = '-7.26%'
    'D7' = '7.666'
    'E7' = '-2.53%'
    'D8' = '3.723'
    'E8' = '-2.24%'
    'E9' = '-3.33%'
    'D10' = '0.1656'
    'E10' = '-6.51%'
    'D11' = '0.07725'
    'E11' = '-1.08%'
    'D12' = '0.08043'
    'E12' = '-7.14%'
    'D13' = '0.03037'
    'E13' = '-3.48%'
    'E14' = '-0.12%'
    'D15' = '0.001507'
    'E15' = '-1.09%'
    'D16' = '0.005862'
    'E16' = '1.49%'
    'D17' = '3.482'
    'E17' = '0.59%'
    'D18' = '2.081'
    'E18' = '-3.41%'
    'D19' = '0.3313'
    'E19' = '0.04%'
    'D20' = '0.1297'
    'E20' = '-1.65%'
    'D21' = '4.043'
    'E21' = '-5.86%'
    'E22' = '5.53%'
    'D23' = '0.04523'
    'E23' = '-0.84%'
    'D24' = '0.001215'
    'E24' = '-0.91%'
    'D25' = '0.004002'
    'E25' = '-9.89%'
    'D26' = '0.0001250'
    'E26' = '-0.04%'
    'D39' = '0.01592'
    'E39' = '-7.01%'
    'D40' = '0.04387'
    'E40' = '-8.13%'
    'D41' = '0.007342'
    'E41' = '-2.39%'
    'D42' = '0.1306'
    'E42' = '-3.57%'
    'D43' = '0.007676'
    'D44' = '0.002014'
    'E44' = '-14.00%'
    'D45' = '0.009201'
    'E45' = '-20.93%'
    'D46' = '0.00005918'
    'E46' = '-5.47%'
    'E47' = '-0.04%'
    'E48' = '172.73%'
    'D49' = '0.003002'
    'E49' = '-3.24%'
    'E50' = '-0.04%'
    'E51' = '-0.04%'
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
}
